$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 14.452709
$ws.Range("H2").Value = 43.358127
$ws.Range("I2").Value = 0.1476906377370901
$ws.Range("J2").Value = 0.1476906377370901
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 45.931316
$ws.Range("N2").Value = 137.793948
$ws.Range("O2").Value = 0.9874217014725413
$ws.Range("P2").Value = 0.9874217014725412
$ws.Range("Q2").Value = 663.831944135044
$ws.Range("R2").Value = 5974.487497215396
$ws.Range("S2").Value = 0.1458329408059222
$ws.Range("T2").Value = 0.1458329408059222

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 14.452709
$ws.Range("H3").Value = 43.358127
$ws.Range("I3").Value = 0.1476906377370901
$ws.Range("J3").Value = 0.1476906377370901
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 0.5850973333333334
$ws.Range("N3").Value = 1.755292
$ws.Range("O3").Value = 0.01257829852745884
$ws.Range("P3").Value = 0.01257829852745884
$ws.Range("Q3").Value = 8.456241495342667
$ws.Range("R3").Value = 76.106173458084
$ws.Range("S3").Value = 0.001857696931167897
$ws.Range("T3").Value = 0.001857696931167898

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 45.91529066666667
$ws.Range("H4").Value = 137.745872
$ws.Range("I4").Value = 0.4692032864180593
$ws.Range("J4").Value = 0.4692032864180593
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 45.931316
$ws.Range("N4").Value = 137.793948
$ws.Range("O4").Value = 0.9874217014725413
$ws.Range("P4").Value = 0.9874217014725412
$ws.Range("Q4").Value = 2108.949724842518
$ws.Range("R4").Value = 18980.54752358266
$ws.Range("S4").Value = 0.4633015074114282
$ws.Range("T4").Value = 0.4633015074114282

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 45.91529066666667
$ws.Range("H5").Value = 137.745872
$ws.Range("I5").Value = 0.4692032864180593
$ws.Range("J5").Value = 0.4692032864180593
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 0.5850973333333334
$ws.Range("N5").Value = 1.755292
$ws.Range("O5").Value = 0.01257829852745884
$ws.Range("P5").Value = 0.01257829852745884
$ws.Range("Q5").Value = 26.86491412829156
$ws.Range("R5").Value = 241.784227154624
$ws.Range("S5").Value = 0.005901779006631122
$ws.Range("T5").Value = 0.005901779006631123

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 36.015061
$ws.Range("H6").Value = 108.045183
$ws.Range("I6").Value = 0.3680339324088102
$ws.Range("J6").Value = 0.3680339324088103
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 45.931316
$ws.Range("N6").Value = 137.793948
$ws.Range("O6").Value = 0.9874217014725413
$ws.Range("P6").Value = 0.9874217014725412
$ws.Range("Q6").Value = 1654.219147550276
$ws.Range("R6").Value = 14887.97232795249
$ws.Range("S6").Value = 0.3634046917387376
$ws.Range("T6").Value = 0.3634046917387376

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 36.015061
$ws.Range("H7").Value = 108.045183
$ws.Range("I7").Value = 0.3680339324088102
$ws.Range("J7").Value = 0.3680339324088103
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 0.5850973333333334
$ws.Range("N7").Value = 1.755292
$ws.Range("O7").Value = 0.01257829852745884
$ws.Range("P7").Value = 0.01257829852745884
$ws.Range("Q7").Value = 21.07231615093734
$ws.Range("R7").Value = 189.650845358436
$ws.Range("S7").Value = 0.004629240670072623
$ws.Range("T7").Value = 0.004629240670072624

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 1.474929666666667
$ws.Range("H8").Value = 4.424789000000001
$ws.Range("I8").Value = 0.01507214343604052
$ws.Range("J8").Value = 0.01507214343604052
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 45.931316
$ws.Range("N8").Value = 137.793948
$ws.Range("O8").Value = 0.9874217014725413
$ws.Range("P8").Value = 0.9874217014725412
$ws.Range("Q8").Value = 67.74546059744135
$ws.Range("R8").Value = 609.7091453769721
$ws.Range("S8").Value = 0.01488256151645333
$ws.Range("T8").Value = 0.01488256151645333

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 1.474929666666667
$ws.Range("H9").Value = 4.424789000000001
$ws.Range("I9").Value = 0.01507214343604052
$ws.Range("J9").Value = 0.01507214343604052
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 0.5850973333333334
$ws.Range("N9").Value = 1.755292
$ws.Range("O9").Value = 0.01257829852745884
$ws.Range("P9").Value = 0.01257829852745884
$ws.Range("Q9").Value = 0.8629774148208891
$ws.Range("R9").Value = 7.766796733388001
$ws.Range("S9").Value = 0.0001895819195871969
$ws.Range("T9").Value = 0.0001895819195871969

